$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 3551
$ws.Range("F5").Value = 3551
$ws.Range("F6").Value = 255
$ws.Range("F7").Value = 5069
$ws.Range("F8").Value = 5069
$ws.Range("F9").Value = 516
$ws.Range("F10").Value = 351
$ws.Range("F11").Value = 197
$ws.Range("F12").Value = 689
$ws.Range("F17").Value = 313
$ws.Range("F23").Value = 4898
$ws.Range("F24").Value = 4898
$ws.Range("F27").Value = 11
$ws.Range("F28").Value = 6022
$ws.Range("F32").Value = 331
$ws.Range("F33").Value = 708
$ws.Range("F34").Value = 4443
$ws.Range("F37").Value = 140
$ws.Range("F42").Value = 865
$ws.Range("F43").Value = 971

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 46
$ws.Range("F3").Value = 25

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1114

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1114
$ws.Range("F7").Value = 3551
$ws.Range("F8").Value = 3551
$ws.Range("F9").Value = 255
$ws.Range("F10").Value = 5069
$ws.Range("F11").Value = 5069
$ws.Range("F12").Value = 516
$ws.Range("F13").Value = 351
$ws.Range("F14").Value = 197
$ws.Range("F15").Value = 689
$ws.Range("F20").Value = 313
$ws.Range("F22").Value = 46
$ws.Range("F27").Value = 4898
$ws.Range("F28").Value = 4898
$ws.Range("F31").Value = 11
$ws.Range("F32").Value = 6022
$ws.Range("F36").Value = 331
$ws.Range("F37").Value = 708
$ws.Range("F38").Value = 4443
$ws.Range("F40").Value = 25
$ws.Range("F42").Value = 140
$ws.Range("F47").Value = 865
$ws.Range("F48").Value = 971
